$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold price text that must stay exact (e.g. "1.10",
# "0.0000177", thousand-dot-separated numbers) -- force text format so
# Excel does not renormalize them into floating point numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.541.23"
$ws.Range("E2").Value = "  +3.72%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.406.95"
$ws.Range("E3").Value = "  +1.00%  "
$ws.Range("E4").Value = "  +0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.36"
$ws.Range("E5").Value = "  +1.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.69"
$ws.Range("E6").Value = "  +4.84%  "
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("E8").Value = "  +1.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.434.35"
$ws.Range("E9").Value = "  +2.14%  "
$ws.Range("E10").Value = "  +4.51%  "
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.22"
$ws.Range("E12").Value = "  +2.20%  "
$ws.Range("E13").Value = "  +3.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.44"
$ws.Range("E14").Value = "  +6.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000177"
$ws.Range("E15").Value = "  +5.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.902.28"
$ws.Range("E16").Value = "  +3.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.320.94"
$ws.Range("E17").Value = "  +3.92%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.431.03"
$ws.Range("E18").Value = "  +2.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.89"
$ws.Range("E19").Value = "  -2.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.94"
$ws.Range("E20").Value = "  +3.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "327.53"
$ws.Range("E21").Value = "  +1.51%  "
$ws.Range("E22").Value = "  +1.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.06"
$ws.Range("E23").Value = "  +11.51%  "
$ws.Range("E24").Value = "  -0.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.58"
$ws.Range("E25").Value = "  +2.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "623.68"
$ws.Range("E26").Value = "  +10.14%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.10"
$ws.Range("E27").Value = "  +10.40%  "
$ws.Range("B28").Value = "BabyDogeCoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₆0554"
$ws.Range("E28").Value = "  +92.25%  "
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.52"
$ws.Range("E29").Value = "  +4.59%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0980"
$ws.Range("E30").Value = "  +5.24%  "
$ws.Range("B31").Value = "WrappedeETH"
$ws.Range("C31").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.561.11"
$ws.Range("E31").Value = "  +2.60%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.18"
$ws.Range("E32").Value = "  +2.32%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.41"
$ws.Range("E33").Value = "  +6.89%  "
$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.85"
$ws.Range("E34").Value = "  +2.93%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.137"
$ws.Range("E35").Value = "  +3.49%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.49"
$ws.Range("E36").Value = "  +2.30%  "
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.996"
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.75"
$ws.Range("E38").Value = "  +3.67%  "
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.374"
$ws.Range("E39").Value = "  +1.56%  "
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "151.92"
$ws.Range("E40").Value = "  -1.28%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.38"
$ws.Range("E41").Value = "  +5.86%  "
$ws.Range("B42").Value = "EthereumClassic"
$ws.Range("C42").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "18.61"
$ws.Range("E42").Value = "  +2.24%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.76"
$ws.Range("E43").Value = "  +12.64%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.75"
$ws.Range("E44").Value = "  +4.78%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "14.74"
$ws.Range("E46").Value = "  +25.88%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "144.29"
$ws.Range("E47").Value = "  +2.77%  "
$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.59"
$ws.Range("E48").Value = "  +1.62%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.49"
$ws.Range("E49").Value = "  +6.43%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.598"
$ws.Range("E50").Value = "  +1.77%  "
$ws.Range("B51").Value = "Hedera"
$ws.Range("C51").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0514"
$ws.Range("E51").Value = "  +2.48%  "
